# Fill in the previously-missing xG_home, xG_away, goals_home and goals_away
# figures for the last six recorded Torino matches (rows 10-15), bringing the
# sheet's D:G columns up to date with the rest of the season's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the target block as text first so the numeric-looking values
# (e.g. "1.63863") are stored as text, matching the rest of the sheet
# (columns D-G throughout the workbook are written as text, not numbers).
$targetRange = $ws.Range("D10:G15")
$targetRange.NumberFormat = "@"

# xG_home (column D) for rows 10-15
$ws.Range("D10").Value = "1.63863"
$ws.Range("D11").Value = "2.52836"
$ws.Range("D12").Value = "1.07631"
$ws.Range("D13").Value = "3.04581"
$ws.Range("D14").Value = "1.10564"
$ws.Range("D15").Value = "0.896358"

# xG_away (column E) for rows 10-15
$ws.Range("E10").Value = "1.38058"
$ws.Range("E11").Value = "1.34238"
$ws.Range("E12").Value = "1.47406"
$ws.Range("E13").Value = "0.975058"
$ws.Range("E14").Value = "2.17731"
$ws.Range("E15").Value = "0.390204"

# goals_home (column F) for rows 10-15
$ws.Range("F10").Value = "2"
$ws.Range("F11").Value = "2"
$ws.Range("F12").Value = "2"
$ws.Range("F13").Value = "3"
$ws.Range("F14").Value = "1"
$ws.Range("F15").Value = "1"

# goals_away (column G) for rows 10-15
$ws.Range("G10").Value = "2"
$ws.Range("G11").Value = "1"
$ws.Range("G12").Value = "3"
$ws.Range("G13").Value = "1"
$ws.Range("G14").Value = "1"
$ws.Range("G15").Value = "1"

# Restore the default (Normal) cell style so the formatting change used only
# to coerce text storage doesn't linger on these cells.
$targetRange.Style = "Normal"
